# Applies the 想去人数 (attendee-count) updates captured in the diff
# for the gh-pages data refresh at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 460
$ws.Range("F5").Value = 67
$ws.Range("F6").Value = 3
$ws.Range("F8").Value = 506
$ws.Range("F9").Value = 106
$ws.Range("F10").Value = 1178
$ws.Range("F11").Value = 178
$ws.Range("F12").Value = 215
$ws.Range("F17").Value = 92
$ws.Range("F18").Value = 241
$ws.Range("F19").Value = 1651
$ws.Range("F20").Value = 609
$ws.Range("F22").Value = 196
$ws.Range("F23").Value = 1921
$ws.Range("F26").Value = 918
$ws.Range("F27").Value = 1203
$ws.Range("F29").Value = 1901
$ws.Range("F31").Value = 1600
$ws.Range("F33").Value = 112
$ws.Range("F34").Value = 637
$ws.Range("F35").Value = 857
$ws.Range("F36").Value = 1766
$ws.Range("F37").Value = 885
$ws.Range("F38").Value = 1782
$ws.Range("F39").Value = 198
$ws.Range("F41").Value = 835
$ws.Range("F42").Value = 37
$ws.Range("F43").Value = 839
$ws.Range("F44").Value = 787
$ws.Range("F45").Value = 995
$ws.Range("F46").Value = 37
$ws.Range("F47").Value = 431
$ws.Range("F48").Value = 3319

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 3
$ws.Range("F12").Value = 790

# Sheet: 全部类型 (All types, combined view)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 460
$ws.Range("F5").Value = 67
$ws.Range("F9").Value = 506
$ws.Range("F10").Value = 106
$ws.Range("F11").Value = 1180
$ws.Range("F12").Value = 178
$ws.Range("F13").Value = 215
$ws.Range("F18").Value = 92
$ws.Range("F20").Value = 241
$ws.Range("F21").Value = 1651
$ws.Range("F22").Value = 609
$ws.Range("F24").Value = 196
$ws.Range("F25").Value = 1921
$ws.Range("F28").Value = 1203
$ws.Range("F30").Value = 1600
$ws.Range("F32").Value = 112
$ws.Range("F33").Value = 790
$ws.Range("F35").Value = 637
$ws.Range("F36").Value = 857
$ws.Range("F37").Value = 1766
$ws.Range("F39").Value = 885
$ws.Range("F40").Value = 1782
$ws.Range("F41").Value = 835
$ws.Range("F42").Value = 839
$ws.Range("F43").Value = 787
$ws.Range("F44").Value = 995
$ws.Range("F45").Value = 431
$ws.Range("F48").Value = 3319
